$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ========================================================================
# 1) Cell values
# ========================================================================

# New job row (#5): Core Java Developer ----------------------------------
$ws.Range("A22").Value = 5
$ws.Range("B22").Value = 'Core Java Developer'
$ws.Range("C22").Value = 'Pune/Bangalore'
$ws.Range("D22").Value = 'Job Title: Senior Engineer – Software Development (Core Java)'

# Header detail lines (plain, default style) -----------------------------
$ws.Range("D24").Value = 'Location: Pune/Bangalore'
$ws.Range("D26").Value = 'Experience: 4+ Years'
$ws.Range("D28").Value = 'Duration: Full time'
$ws.Range("D32").Value = 'Job Description:'

# Job description intro paragraph ----------------------------------------
$ws.Range("D38").Value = 'Our clients include banks, telecom, IT, and mid- sized Enterprises across the globe. We are predominantly working with Banking clients Top Private & Public Sector Banks in India and Overseas (Middle east, Africa, South east regions). We have over 50 happy customers and we plan to increase our customer base to 500+ in the next 24 to 36 months.'

# Responsibilities bullet list (D41 stays blank, like the source row) ----
$ws.Range("D42").Value = 'Design, develop, and maintain Java-based applications, ensuring optimal performance, reliability, and scalability.'
$ws.Range("D43").Value = 'Write clean, efficient, and well-documented code following industry best practices and coding standards.'
$ws.Range("D44").Value = 'Participate in the entire software development lifecycle, including requirements analysis, design, implementation, testing, deployment, and support.'
$ws.Range("D45").Value = 'Collaborate with product managers, business analysts, and other stakeholders to understand requirements and translate them into technical solutions.'
$ws.Range("D46").Value = 'Conduct code reviews, provide constructive feedback, and mentor junior team members to promote continuous improvement and knowledge sharing.'
$ws.Range("D47").Value = 'Troubleshoot and debug issues reported by clients or detected during testing, and implement timely and effective solutions.'
$ws.Range("D48").Value = 'Work closely with QA engineers to ensure the quality of software deliverables through thorough testing and validation.'
$ws.Range("D49").Value = 'Contribute to architectural design discussions and decisions and participate in team technical discussions.'
$ws.Range("D50").Value = 'Collaborate with DevOps engineers to automate deployment processes and enhance system monitoring and performance optimization.'

# "Skill Requirement" heading ---------------------------------------------
$ws.Range("D52").Value = 'Skill Requirement'

# Skill bullet list (D53 stays blank, like the source row) ---------------
$ws.Range("D54").Value = 'The ideal candidate should have passion for building products, solving problems, and building data pipeline.'
$ws.Range("D55").Value = 'Proficiency in version 8 and higher Java.'
$ws.Range("D56").Value = 'Experience in Clojure, Scala or Java, knowledge of Spark, Flink.'
$ws.Range("D57").Value = 'The basics must be very strong - design, coding, testing, and debugging skills.'
$ws.Range("D58").Value = 'Proficiency in web application development using Java-based technologies (Servlets, JSP, etc.).'
$ws.Range("D59").Value = 'Familiarity with relational databases (e.g., MySQL, PostgreSQL) and proficiency in SQL.'
$ws.Range("D60").Value = 'BS in CS/EE/CE or related field from a top institution'
$ws.Range("D61").Value = '4+ years hands-on experience in Java, data structures and algorithms on Linux'
$ws.Range("D62").Value = 'Experience/knowledge with Microservices, Docker, Kubernetes, agile methodologies and tools (e.g., Scrum, JIRA) experience is a plus'
$ws.Range("D63").Value = 'Familiarity with cloud platforms (e.g., AWS, Azure) and microservices architecture is desirable.'
$ws.Range("D64").Value = 'A demonstrable understanding of software development concepts, problem break down, project management, and good communication.'
$ws.Range("D65").Value = 'Experience will product build life cycle of developing, debugging, optimizing and maintaining code. '

# ========================================================================
# 2) Cell formatting - build each new style once, then fan it out via
#    copy / paste-special (formats only) so no orphan styles are left
#    behind in the shared style table.
# ========================================================================

# Style A: wrap + vertical-center (intro paragraph + blank spacer rows) --
$ws.Range("D38").WrapText = $true
$ws.Range("D38").VerticalAlignment = -4108
$ws.Range("D38").Copy()
$ws.Range("D39:D40").PasteSpecial(-4122)
$ws.Range("D51").PasteSpecial(-4122)
$ws.Range("D66").PasteSpecial(-4122)

# Style B: wrap + vertical-center + left-indent (bullet lists) -----------
$ws.Range("D41").WrapText = $true
$ws.Range("D41").VerticalAlignment = -4108
$ws.Range("D41").IndentLevel = 1
$ws.Range("D41").HorizontalAlignment = -4131
$ws.Range("D41").Copy()
$ws.Range("D42:D50").PasteSpecial(-4122)
$ws.Range("D53:D65").PasteSpecial(-4122)

# Style C: bold dark-blue font + wrap + vertical-center ("Skill Requirement")
$ws.Range("D52").WrapText = $true
$ws.Range("D52").VerticalAlignment = -4108
$ws.Range("D52").Font.Bold = $true
$ws.Range("D52").Font.Color = 10703646

# ========================================================================
# 3) Row heights for the wrapped paragraph rows
# ========================================================================

$ws.Rows.Item(38).RowHeight = 43.2
$ws.Rows.Item(44).RowHeight = 28.8
$ws.Rows.Item(45).RowHeight = 28.8
$ws.Rows.Item(46).RowHeight = 28.8

# ========================================================================
# 4) Selection, matching the saved worksheet view
# ========================================================================

$ws.Range("D38").Select()

